{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces each of the 100 arithmetic-answer cells in the single table\n// with its new expression, in row-major order (20 rows x 5 cols),\n// matching the authoritative diff.\nconst newValues = [\n  \"95-6=89\",\n  \"24-6=18\",\n  \"71-3=68\",\n  \"24+48=72\",\n  \"62-48=14\",\n  \"36+36=72\",\n  \"39+34=73\",\n  \"16+59=75\",\n  \"98-49=49\",\n  \"83-4=79\",\n  \"46-18=28\",\n  \"58+18=76\",\n  \"90-31=59\",\n  \"83-36=47\",\n  \"81-22=59\",\n  \"47+39=86\",\n  \"8+48=56\",\n  \"63+28=91\",\n  \"24+7=31\",\n  \"29+54=83\",\n  \"71-39=32\",\n  \"64-27=37\",\n  \"16+79=95\",\n  \"60-5=55\",\n  \"92-9=83\",\n  \"73-67=6\",\n  \"17+17=34\",\n  \"87-49=38\",\n  \"59+28=87\",\n  \"83-25=58\",\n  \"93-89=4\",\n  \"81-22=59\",\n  \"93-45=48\",\n  \"11-8=3\",\n  \"33+29=62\",\n  \"62-49=13\",\n  \"64-35=29\",\n  \"45+6=51\",\n  \"8+17=25\",\n  \"36+45=81\",\n  \"42-37=5\",\n  \"55-16=39\",\n  \"27+49=76\",\n  \"52-16=36\",\n  \"23+39=62\",\n  \"8+85=93\",\n  \"61-15=46\",\n  \"56+38=94\",\n  \"81-57=24\",\n  \"38+29=67\",\n  \"61-43=18\",\n  \"50-3=47\",\n  \"55+9=64\",\n  \"69+27=96\",\n  \"42-4=38\",\n  \"94-76=18\",\n  \"69+22=91\",\n  \"22-19=3\",\n  \"9+89=98\",\n  \"90-13=77\",\n  \"83+8=91\",\n  \"96-39=57\",\n  \"55+37=92\",\n  \"39+7=46\",\n  \"37+34=71\",\n  \"6+28=34\",\n  \"57-38=19\",\n  \"62-27=35\",\n  \"10-8=2\",\n  \"62-54=8\",\n  \"63+8=71\",\n  \"61-13=48\",\n  \"70-62=8\",\n  \"58+3=61\",\n  \"27-8=19\",\n  \"91-19=72\",\n  \"71-34=37\",\n  \"8+25=33\",\n  \"14+78=92\",\n  \"7+16=23\",\n  \"85-67=18\",\n  \"27+4=31\",\n  \"45-38=7\",\n  \"17+74=91\",\n  \"19+5=24\",\n  \"67-19=48\",\n  \"91-74=17\",\n  \"51-15=36\",\n  \"3+59=62\",\n  \"90-65=25\",\n  \"46+46=92\",\n  \"44-18=26\",\n  \"19+75=94\",\n  \"37+34=71\",\n  \"7+27=34\",\n  \"68-49=19\",\n  \"61-18=43\",\n  \"13+68=81\",\n  \"81-5=76\",\n  \"7+54=61\"\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst ROWS = table.rowCount;\nconst COLS = 5;\nif (ROWS * COLS !== newValues.length) {\n  throw new Error(`Table shape ${ROWS}x${COLS} does not match ${newValues.length} replacement values.`);\n}\n\nlet idx = 0;\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    table.getCell(r, c).value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# PowerShell (Word COM interop) edit script.\n# Replaces each of the 100 arithmetic-answer cells in the single table\n# with its new expression, in row-major order (20 rows x 5 cols),\n# matching the authoritative diff.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    '95-6=89',\n    '24-6=18',\n    '71-3=68',\n    '24+48=72',\n    '62-48=14',\n    '36+36=72',\n    '39+34=73',\n    '16+59=75',\n    '98-49=49',\n    '83-4=79',\n    '46-18=28',\n    '58+18=76',\n    '90-31=59',\n    '83-36=47',\n    '81-22=59',\n    '47+39=86',\n    '8+48=56',\n    '63+28=91',\n    '24+7=31',\n    '29+54=83',\n    '71-39=32',\n    '64-27=37',\n    '16+79=95',\n    '60-5=55',\n    '92-9=83',\n    '73-67=6',\n    '17+17=34',\n    '87-49=38',\n    '59+28=87',\n    '83-25=58',\n    '93-89=4',\n    '81-22=59',\n    '93-45=48',\n    '11-8=3',\n    '33+29=62',\n    '62-49=13',\n    '64-35=29',\n    '45+6=51',\n    '8+17=25',\n    '36+45=81',\n    '42-37=5',\n    '55-16=39',\n    '27+49=76',\n    '52-16=36',\n    '23+39=62',\n    '8+85=93',\n    '61-15=46',\n    '56+38=94',\n    '81-57=24',\n    '38+29=67',\n    '61-43=18',\n    '50-3=47',\n    '55+9=64',\n    '69+27=96',\n    '42-4=38',\n    '94-76=18',\n    '69+22=91',\n    '22-19=3',\n    '9+89=98',\n    '90-13=77',\n    '83+8=91',\n    '96-39=57',\n    '55+37=92',\n    '39+7=46',\n    '37+34=71',\n    '6+28=34',\n    '57-38=19',\n    '62-27=35',\n    '10-8=2',\n    '62-54=8',\n    '63+8=71',\n    '61-13=48',\n    '70-62=8',\n    '58+3=61',\n    '27-8=19',\n    '91-19=72',\n    '71-34=37',\n    '8+25=33',\n    '14+78=92',\n    '7+16=23',\n    '85-67=18',\n    '27+4=31',\n    '45-38=7',\n    '17+74=91',\n    '19+5=24',\n    '67-19=48',\n    '91-74=17',\n    '51-15=36',\n    '3+59=62',\n    '90-65=25',\n    '46+46=92',\n    '44-18=26',\n    '19+75=94',\n    '37+34=71',\n    '7+27=34',\n    '68-49=19',\n    '61-18=43',\n    '13+68=81',\n    '81-5=76',\n    '7+54=61'\n)\n\n$cols = 5\nfor ($i = 0; $i -lt $newValues.Length; $i++) {\n    $row = [int]([math]::Floor($i / $cols)) + 1\n    $col = ($i % $cols) + 1\n    $table.Cell($row, $col).Range.Text = $newValues[$i]\n}\n\n"}
